$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text values (avoid Excel auto-converting numeric-looking strings to numbers)
# by prefixing with a literal apostrophe (quote-prefix) and then resetting the cell
# style back to "Normal" so no extra formatting/style is introduced.

$ws.Range('D2').Value = "'41.531.38"
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = "'  +0.46%  "
$ws.Range('E2').Style = 'Normal'

$ws.Range('D3').Value = "'2.495.45"
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = "'  +1.12%  "
$ws.Range('E3').Style = 'Normal'

$ws.Range('D4').Value = "'0.998"
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = "'  -0.20%  "
$ws.Range('E4').Style = 'Normal'

$ws.Range('D5').Value = "'313.66"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = "'  +0.27%  "
$ws.Range('E5').Style = 'Normal'

$ws.Range('D6').Value = "'93.64"
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = "'  -1.02%  "
$ws.Range('E6').Style = 'Normal'

$ws.Range('E7').Value = "'  -0.96%  "
$ws.Range('E7').Style = 'Normal'

$ws.Range('E8').Value = "'  -0.18%  "
$ws.Range('E8').Style = 'Normal'

$ws.Range('D9').Value = "'0.501"
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = "'  -0.41%  "
$ws.Range('E9').Style = 'Normal'

$ws.Range('D10').Value = "'32.73"
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = "'  -2.72%  "
$ws.Range('E10').Style = 'Normal'

$ws.Range('D11').Value = "'0.0786"
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = "'  +0.73%  "
$ws.Range('E11').Style = 'Normal'

$ws.Range('E12').Value = "'  +1.86%  "
$ws.Range('E12').Style = 'Normal'

$ws.Range('D13').Value = "'2.873.31"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = "'  +0.88%  "
$ws.Range('E13').Style = 'Normal'

$ws.Range('D14').Value = "'6.89"
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = "'  -2.19%  "
$ws.Range('E14').Style = 'Normal'

$ws.Range('D15').Value = "'16.08"
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = "'  +9.76%  "
$ws.Range('E15').Style = 'Normal'

$ws.Range('D16').Value = "'2.490.11"
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = "'  +1.02%  "
$ws.Range('E16').Style = 'Normal'

$ws.Range('D17').Value = "'0.759"
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = "'  -3.56%  "
$ws.Range('E17').Style = 'Normal'

$ws.Range('D18').Value = "'41.582.05"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = "'  +0.71%  "
$ws.Range('E18').Style = 'Normal'

$ws.Range('D19').Value = "'6.36"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = "'  +0.26%  "
$ws.Range('E19').Style = 'Normal'

$ws.Range('D20').Value = "'0.0₃0931"
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = "'  +1.01%  "
$ws.Range('E20').Style = 'Normal'

$ws.Range('D21').Value = "'71.36"
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = "'  +4.76%  "
$ws.Range('E21').Style = 'Normal'

$ws.Range('E22').Value = "'  -2.74%  "
$ws.Range('E22').Style = 'Normal'

$ws.Range('D23').Value = "'236.26"
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = "'  -0.55%  "
$ws.Range('E23').Style = 'Normal'

$ws.Range('D24').Value = "'2.72"
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = "'  -2.67%  "
$ws.Range('E24').Style = 'Normal'

$ws.Range('B25').Value = "'ImmutableX"
$ws.Range('B25').Style = 'Normal'
$ws.Range('C25').Value = "'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range('C25').Style = 'Normal'
$ws.Range('D25').Value = "'1.93"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = "'  -0.36%  "
$ws.Range('E25').Style = 'Normal'

$ws.Range('B26').Value = "'Dai"
$ws.Range('B26').Style = 'Normal'
$ws.Range('C26').Value = "'https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range('C26').Style = 'Normal'
$ws.Range('D26').Value = "'1.00"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = "'  -0.07%  "
$ws.Range('E26').Style = 'Normal'

$ws.Range('D27').Value = "'25.37"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = "'  +3.66%  "
$ws.Range('E27').Style = 'Normal'

$ws.Range('E28').Value = "'  -0.55%  "
$ws.Range('E28').Style = 'Normal'

$ws.Range('D29').Value = "'9.70"
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = "'  -0.14%  "
$ws.Range('E29').Style = 'Normal'

$ws.Range('D30').Value = "'36.29"
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = "'  +0.53%  "
$ws.Range('E30').Style = 'Normal'

$ws.Range('D31').Value = "'157.75"
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = "'  +2.69%  "
$ws.Range('E31').Style = 'Normal'

$ws.Range('D32').Value = "'5.46"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = "'  -2.05%  "
$ws.Range('E32').Style = 'Normal'

$ws.Range('D33').Value = "'2.58"
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = "'  -1.02%  "
$ws.Range('E33').Style = 'Normal'

$ws.Range('D34').Value = "'0.0760"
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = "'  +0.42%  "
$ws.Range('E34').Style = 'Normal'

$ws.Range('D35').Value = "'18.02"
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = "'  +5.98%  "
$ws.Range('E35').Style = 'Normal'

$ws.Range('E36').Value = "'  -6.19%  "
$ws.Range('E36').Style = 'Normal'

$ws.Range('D37').Value = "'2.96"
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = "'  -1.94%  "
$ws.Range('E37').Style = 'Normal'

$ws.Range('B38').Value = "'ARBITRUM"
$ws.Range('B38').Style = 'Normal'
$ws.Range('C38').Value = "'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range('C38').Style = 'Normal'
$ws.Range('D38').Value = "'1.86"
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = "'  -1.80%  "
$ws.Range('E38').Style = 'Normal'

$ws.Range('B39').Value = "'Kaspa"
$ws.Range('B39').Style = 'Normal'
$ws.Range('C39').Value = "'https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range('C39').Style = 'Normal'
$ws.Range('D39').Value = "'0.105"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = "'  +1.40%  "
$ws.Range('E39').Style = 'Normal'

$ws.Range('E40').Value = "'  -0.12%  "
$ws.Range('E40').Style = 'Normal'

$ws.Range('D41').Value = "'4.14"
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = "'  -4.35%  "
$ws.Range('E41').Style = 'Normal'

$ws.Range('E42').Value = "'  -0.26%  "
$ws.Range('E42').Style = 'Normal'

$ws.Range('D43').Value = "'20.03"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = "'  -6.09%  "
$ws.Range('E43').Style = 'Normal'

$ws.Range('D44').Value = "'1.974.25"
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = "'  -0.67%  "
$ws.Range('E44').Style = 'Normal'

$ws.Range('D45').Value = "'0.0285"
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = "'  -0.19%  "
$ws.Range('E45').Style = 'Normal'

$ws.Range('D46').Value = "'3.00"
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = "'  -2.77%  "
$ws.Range('E46').Style = 'Normal'

$ws.Range('D47').Value = "'8.89"
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = "'  +1.95%  "
$ws.Range('E47').Style = 'Normal'

$ws.Range('D48').Value = "'2.730.31"
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = "'  +0.83%  "
$ws.Range('E48').Style = 'Normal'

$ws.Range('D49').Value = "'96.83"
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = "'  -0.36%  "
$ws.Range('E49').Style = 'Normal'

$ws.Range('D50').Value = "'68.19"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = "'  -2.60%  "
$ws.Range('E50').Style = 'Normal'

$ws.Range('D51').Value = "'73.84"
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = "'  -3.06%  "
$ws.Range('E51').Style = 'Normal'
